# Commit: "Fruta / hortaliza, semanal"
# The underlying data rows (2-38) get re-shuffled: each row's varying
# attributes (date, variety, quality, volume, prices, unit, price/kg, kg-or-units)
# are replaced by the values that used to live on a different row, per a
# fixed permutation. Columns A,B,C,E,F,G,O,R are constant across all rows
# for this market/product so they are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row (values read from source row in the
# ORIGINAL workbook get written to the destination row)
$map = @{
    2=11; 3=2; 4=32; 5=9; 6=17; 7=36; 8=37; 9=19; 10=12;
    11=35; 12=22; 13=33; 14=38; 15=30; 16=23; 17=20; 18=10; 19=16; 20=14;
    21=21; 22=27; 23=28; 24=24; 25=6; 26=15; 27=7; 28=3; 29=26;
    30=5; 31=13; 32=31; 33=8; 34=34; 35=29; 36=4; 37=18; 38=25
}

# Columns whose values move together with the row (by column index)
$cols = @(4, 8, 9, 10, 11, 12, 13, 14, 16, 17)

# Snapshot every source cell's current value before any writes happen,
# so that writes to one destination never clobber a value still needed
# as a source for another destination (the mapping is a permutation).
$snapshot = @{}
foreach ($r in 2..38) {
    foreach ($c in $cols) {
        $key = "$r,$c"
        $snapshot[$key] = $ws.Cells.Item($r, $c).Value2
    }
}

foreach ($destRow in 2..38) {
    $srcRow = $map[$destRow]
    foreach ($c in $cols) {
        $val = $snapshot["$srcRow,$c"]
        $ws.Cells.Item($destRow, $c).Value = $val
    }
}

Write-Host "done"
